$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.171.56"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.584.45"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "211.31"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "19.20"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "1.806.50"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.569.33"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "63.98"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "26.140.39"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "213.52"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "8.94"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "143.77"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "1.346.82"
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "0.579"
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "0.934"
$ws.Range("E42").Value = "  -16.71%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "1.719.22"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "60.77"
$ws.Range("E46").Value = "  -3.21%  "
$ws.Range("D47").Value = "86.06"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").Value = "  +5.48%  "
$ws.Range("D49").Value = "1.48"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "0.0981"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -1.14%  "
